$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.974.98'
$ws.Range("E2").Value = '  +2.32%  '
$ws.Range("D3").Value = '1.707.76'
$ws.Range("E3").Value = '  +1.54%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3954'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.87%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4037'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.90%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.486'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '52.78'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9999'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08818'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '26.01'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.491'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001357'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.999'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("D17").Value = '1.718.91'
$ws.Range("E17").Value = '  +2.92%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '96.28'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07181'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.365'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9994'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.53%  '
$ws.Range("D24").Value = '24.965.74'
$ws.Range("E24").Value = '  +2.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.980'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.352'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.59%  '
$ws.Range("E27").Value = '  +5.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.220'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +16.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '161.65'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '150.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +9.27%  '
$ws.Range("E31").Value = '  -2.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.581'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +31.07%  '
$ws.Range("D33").Value = '1.902.11'
$ws.Range("E33").Value = '  +2.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08552'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.45%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.048'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.03134'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.183'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2859'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09525'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.84'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8277'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.03'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.484'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.50'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.686'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7402'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.246'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.378'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08739'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +8.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.000'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '139.16'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.20%  '
